# "Finished looking at increasing part"
#
# On the "Alt1" sheet, a new column is inserted before column R (the old
# column R — the FILTER-of-positive-run-starts spill — slides over to
# become column S). The freed-up column R now holds a plain AVERAGE of a
# new helper spill in column Q, which re-filters the FREQUENCY buckets in
# column P down to the ones greater than zero (incremented by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Alt1")
$ws.Activate()

# Insert a new column at R; everything that used to live in R (and to its
# right) shifts one column to the right, becoming column S. Formulas with
# relative/array refs (including the dynamic-array spill anchored at R27)
# are shifted automatically.
$ws.Range("R1").EntireColumn.Insert()

# New helper spill: for each FREQUENCY bucket in P27# that is > 0, add 1.
$ws.Range("Q27").Formula2 = "=1+FILTER(P27#,P27#>0)"

# New plain formula: average of the helper spill above.
$ws.Range("R27").Formula2 = "=AVERAGE(Q27#)"

# Match the author's final selection/scroll state on the sheet.
$ws.Range("G21").Select()
